$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new phone-number values in column C (rows 2 and 3).
# This introduces two new shared strings and two new cells, matching
# the uploaded workbook.
$ws.Range("C2").Value = "0414 270 2375"
$ws.Range("C3").Value = "0412 270 2374"

# Column B was widened slightly so the "nombre destinatario" values
# remain comfortably readable next to the new column.
$ws.Columns("B").ColumnWidth = 24

# Leave the active selection on C3, the last cell touched.
$ws.Range("C3").Select()
